$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.957.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.15%  "

$ws.Range("D3").Value = "'1.859.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.60%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'306.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.24%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.5107"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.31%  "

$ws.Range("D8").Value = "'0.3733"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.41%  "

$ws.Range("D9").Value = "'0.07100"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.40%  "

$ws.Range("D10").Value = "'0.8857"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.94%  "

$ws.Range("D11").Value = "'20.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.17%  "

$ws.Range("D12").Value = "'0.07550"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.15%  "

$ws.Range("D13").Value = "'1.853.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.82%  "

$ws.Range("D14").Value = "'5.299"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.10%  "

$ws.Range("D15").Value = "'88.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.89%  "

$ws.Range("D16").Value = "'1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D17").Value = "'0.000008416"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.46%  "

$ws.Range("D18").Value = "'14.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.55%  "

$ws.Range("E19").Value = "  +0.34%  "

$ws.Range("D20").Value = "'27.011.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.01%  "

$ws.Range("D21").Value = "'5.046"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.38%  "

$ws.Range("D22").Value = "'2.102.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.07%  "

$ws.Range("D23").Value = "'10.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.99%  "

$ws.Range("D24").Value = "'6.467"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.04%  "

$ws.Range("D25").Value = "'149.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.41%  "

$ws.Range("D26").Value = "'1.838"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "

$ws.Range("D27").Value = "'17.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.30%  "

$ws.Range("D28").Value = "'2.097"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.39%  "

$ws.Range("D29").Value = "'112.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.06%  "

$ws.Range("D30").Value = "'4.673"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.13%  "

$ws.Range("D31").Value = "'4.648"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.22%  "

$ws.Range("D32").Value = "'0.09023"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.53%  "

$ws.Range("D33").Value = "'0.05112"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.47%  "

$ws.Range("D34").Value = "'3.071"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.83%  "

$ws.Range("D35").Value = "'1.151"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.99%  "

$ws.Range("D36").Value = "'0.7314"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.29%  "

$ws.Range("D37").Value = "'0.02045"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.61%  "

$ws.Range("D38").Value = "'2.493"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.99%  "

$ws.Range("D39").Value = "'3.047"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("D40").Value = "'1.071"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.02%  "

$ws.Range("D41").Value = "'0.5308"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.79%  "

$ws.Range("D42").Value = "'6.576"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.81%  "

$ws.Range("D43").Value = "'115.44"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "'8.269"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.42%  "

$ws.Range("D45").Value = "'0.1468"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.99%  "

$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").Value = "'0.4609"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.84%  "

$ws.Range("E48").Value = "  -4.74%  "

$ws.Range("D49").Value = "'1.561"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.50%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'36.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.98%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'64.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.79%  "

